# Re-assessment of detection-tool performance: update the "Detection Tool
# Performance Value" (column C) raw scores for the three criteria in both
# the input table (C5:C7) and the constructed TOPSIS matrix (C29:C31).
# Everything else on the sheet is a formula and recalculates automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 2. Define relevant selection criteria / input scores (rows 5-7) ---
$ws.Range("C5").Value = 0.85
$ws.Range("C6").Value = 0.17
$ws.Range("C7").Value = 0.1

# --- 4. Construct the Matrix (rows 29-31) - same re-assessed values ---
$ws.Range("C29").Value = 0.85
$ws.Range("C30").Value = 0.17
$ws.Range("C31").Value = 0.1

# Re-enter the "copy down" formulas that reference the weighted-matrix rows
# (8. ideal positive / negative separation blocks) so they stay in sync.
$ws.Range("B61").Formula = "=B54"
$ws.Range("C61").Formula = "=C54"
$ws.Range("B62").Formula = "=B55"
$ws.Range("C62").Formula = "=C55"
$ws.Range("B63").Formula = "=B56"
$ws.Range("C63").Formula = "=C56"

$ws.Range("B69").Formula = "=B54"
$ws.Range("C69").Formula = "=C54"
$ws.Range("B70").Formula = "=B55"
$ws.Range("C70").Formula = "=C55"
$ws.Range("B71").Formula = "=B56"
$ws.Range("C71").Formula = "=C56"

# Recalculate so every cached <v> reflects the re-assessment.
$excel.Calculate()

# --- Restore the view: scrolled down to row 97, selection on C32 ---
$win = $excel.ActiveWindow
$win.ScrollRow = 97
$win.ScrollColumn = 1
[void]$ws.Range("C32").Select()
